# Apply updated cryptocurrency data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.453.29'
$ws.Range('E2').Value = '  +1.22%  '

$ws.Range('D3').Value = '1.825.86'
$ws.Range('E3').Value = '  +1.83%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.36'
$ws.Range('E5').Value = '  -0.14%  '

$ws.Range('E6').Value = '  +0.17%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5411'
$ws.Range('E7').Value = '  +0.93%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4037'
$ws.Range('E8').Value = '  +7.12%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07659'
$ws.Range('E9').Value = '  +2.56%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.88'
$ws.Range('E10').Value = '  +0.18%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.114'
$ws.Range('E11').Value = '  +1.41%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.325'
$ws.Range('E12').Value = '  +3.46%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.654'
$ws.Range('E13').Value = '  +5.79%  '

$ws.Range('E14').Value = '  +0.30%  '

$ws.Range('E15').Value = '  +1.13%  '

$ws.Range('D16').Value = '1.828.23'
$ws.Range('E16').Value = '  +2.77%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '89.68'
$ws.Range('E17').Value = '  +0.90%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001078'
$ws.Range('E18').Value = '  +2.00%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06594'
$ws.Range('E19').Value = '  +2.12%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').Value = '  +1.73%  '

$ws.Range('E21').Value = '  +0.17%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.071'
$ws.Range('E22').Value = '  +2.82%  '

$ws.Range('D23').Value = '28.458.56'
$ws.Range('E23').Value = '  +1.24%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.08'

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.224'
$ws.Range('E25').Value = '  +5.95%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.464'
$ws.Range('E26').Value = '  +7.48%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.73'
$ws.Range('E27').Value = '  +1.95%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '157.29'
$ws.Range('E28').Value = '  +1.32%  '

$ws.Range('D29').Value = '2.038.13'

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.90'
$ws.Range('E30').Value = '  +3.08%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1118'
$ws.Range('E31').Value = '  +6.60%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.124'
$ws.Range('E32').Value = '  +0.66%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.677'
$ws.Range('E33').Value = '  +2.02%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07379'
$ws.Range('E34').Value = '  +14.09%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.641'
$ws.Range('E35').Value = '  -0.03%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2245'
$ws.Range('E36').Value = '  -0.75%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02340'
$ws.Range('E37').Value = '  +2.18%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.203'
$ws.Range('E38').Value = '  +3.52%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.842'
$ws.Range('E39').Value = '  +4.20%  '

$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '11.33'
$ws.Range('E40').Value = '  +1.95%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6270'
$ws.Range('E41').Value = '  +1.60%  '

$ws.Range('E42').Value = '  +0.66%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.392'
$ws.Range('E44').Value = '  -3.87%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  +2.02%  '

$ws.Range('E46').Value = '  +0.70%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5849'
$ws.Range('E47').Value = '  +1.17%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.90'
$ws.Range('E48').Value = '  -1.67%  '

$ws.Range('E49').Value = '  +3.87%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.202'
$ws.Range('E50').Value = '  +0.85%  '

$ws.Range('E51').Value = '  +1.17%  '
